# Applies the cryptos.xlsx data refresh described in the commit message
# 'Updated cryptos list on Sun Dec 10 13:50:03 UTC 2023 with GitHub Actions'.
# Only the Price (D) / Volume(1h) (E) text values change, plus the
# Cosmos/Toncoin row swap (rows 28-29, columns B/C/D/E) -- no rows/
# columns are inserted or removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while keeping it stored as text.
# Source cells are inline/shared strings (e.g. '43.937.37', '239.36',
# '0.667'); several of the replacement values parse as plain numbers,
# so Excel would otherwise silently convert them to numeric cells and
# introduce float rounding (239.36 -> 239.36000000000001). Forcing the
# number format to Text ('@') before the assignment keeps them as the
# exact original string whenever that would happen.
function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $text
}

Set-TextCell 'D2' '43.937.37'
Set-TextCell 'E2' '  -0.28%  '
Set-TextCell 'D3' '2.348.95'
Set-TextCell 'E3' '  -0.35%  '
Set-TextCell 'E4' '  +0.10%  '
Set-TextCell 'D5' '239.36'
Set-TextCell 'E5' '  -0.83%  '
Set-TextCell 'D6' '0.667'
Set-TextCell 'E6' '  -3.97%  '
Set-TextCell 'D7' '72.88'
Set-TextCell 'E7' '  -4.64%  '
Set-TextCell 'E8' '  +0.03%  '
Set-TextCell 'E9' '  -4.87%  '
Set-TextCell 'E10' '  -1.23%  '
Set-TextCell 'D11' '59.47'
Set-TextCell 'E11' '  +3.68%  '
Set-TextCell 'D12' '32.89'
Set-TextCell 'E12' '  -0.94%  '
Set-TextCell 'E13' '  -0.31%  '
Set-TextCell 'E14' '  -3.48%  '
Set-TextCell 'D15' '2.693.73'
Set-TextCell 'E15' '  -0.44%  '
Set-TextCell 'D16' '16.11'
Set-TextCell 'E16' '  -4.03%  '
Set-TextCell 'D17' '0.904'
Set-TextCell 'E17' '  -2.59%  '
Set-TextCell 'D18' '2.346.42'
Set-TextCell 'E18' '  -0.42%  '
Set-TextCell 'D19' '43.784.83'
Set-TextCell 'E19' '  -0.27%  '
Set-TextCell 'E20' '  -0.08%  '
Set-TextCell 'E21' '  -0.29%  '
Set-TextCell 'D22' '78.40'
Set-TextCell 'E22' '  +0.71%  '
Set-TextCell 'D23' '251.82'
Set-TextCell 'E23' '  -4.15%  '
Set-TextCell 'E24' '  +0.16%  '
Set-TextCell 'D25' '3.75'
Set-TextCell 'E25' '  +3.13%  '
Set-TextCell 'E26' '  +1.49%  '
Set-TextCell 'D27' '2.49'
Set-TextCell 'E27' '  -1.51%  '
Set-TextCell 'B28' 'Toncoin'
Set-TextCell 'C28' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D28' '2.33'
Set-TextCell 'E28' '  +1.08%  '
Set-TextCell 'B29' 'Cosmos'
Set-TextCell 'C29' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D29' '10.40'
Set-TextCell 'E29' '  -4.89%  '
Set-TextCell 'D30' '176.64'
Set-TextCell 'E30' '  +0.74%  '
Set-TextCell 'E31' '  -4.40%  '
Set-TextCell 'E32' '  -1.45%  '
Set-TextCell 'E33' '  -2.63%  '
Set-TextCell 'D34' '0.0745'
Set-TextCell 'E34' '  -2.47%  '
Set-TextCell 'D35' '5.08'
Set-TextCell 'E35' '  -6.06%  '
Set-TextCell 'D36' '5.34'
Set-TextCell 'E36' '  -2.47%  '
Set-TextCell 'D37' '3.76'
Set-TextCell 'E37' '  -1.69%  '
Set-TextCell 'D38' '6.40'
Set-TextCell 'E38' '  -0.21%  '
Set-TextCell 'E39' '  -2.55%  '
Set-TextCell 'D40' '5.71'
Set-TextCell 'E40' '  +20.38%  '
Set-TextCell 'D41' '0.0271'
Set-TextCell 'E41' '  -4.37%  '
Set-TextCell 'D42' '65.34'
Set-TextCell 'E42' '  +15.61%  '
Set-TextCell 'D43' '9.23'
Set-TextCell 'E43' '  +0.84%  '
Set-TextCell 'E44' '  -2.72%  '
Set-TextCell 'D45' '18.79'
Set-TextCell 'E45' '  -2.83%  '
Set-TextCell 'D46' '0.195'
Set-TextCell 'E46' '  -10.62%  '
Set-TextCell 'E47' '  -0.07%  '
Set-TextCell 'E48' '  -2.96%  '
Set-TextCell 'E49' '  -3.18%  '
Set-TextCell 'D50' '98.12'
Set-TextCell 'E50' '  -3.85%  '
Set-TextCell 'E51' '  -5.03%  '
